$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 4.8
$ws.Range("H2").Value = 1.92
$ws.Range("I2").Value = 2.14
$ws.Range("J2").Value = 3.35
$ws.Range("K2").Value = 4
$ws.Range("P2").Value = 1.86
$ws.Range("Q2").Value = 1.96

# Row 3
$ws.Range("F3").Value = 3.45
$ws.Range("G3").Value = 4.4
$ws.Range("H3").Value = 2.06
$ws.Range("I3").Value = 2.32
$ws.Range("J3").Value = 3.35
$ws.Range("K3").Value = 3.95
$ws.Range("P3").Value = 1.88
$ws.Range("Q3").Value = 1.9

# Row 4
$ws.Range("O4").Value = 1.33
$ws.Range("T4").Value = 1.79
$ws.Range("AH4").Value = 18.5

# Row 5
$ws.Range("F5").Value = 2.38
$ws.Range("G5").Value = 2.42
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 3.4

# Row 6
$ws.Range("F6").Value = 4.9
$ws.Range("I6").Value = 1.9
$ws.Range("J6").Value = 3.55
$ws.Range("M6").Value = 1.09
$ws.Range("P6").Value = 1.78
$ws.Range("Q6").Value = 2.14
$ws.Range("S6").Value = 4.1
$ws.Range("X6").Value = 12.5
$ws.Range("AA6").Value = 22
$ws.Range("AE6").Value = 24
$ws.Range("AG6").Value = 21
$ws.Range("AI6").Value = 50
$ws.Range("AO6").Value = 17

# Row 7
$ws.Range("H7").Value = 3

# Row 8
$ws.Range("U8").Value = 1.94

# Row 10
$ws.Range("U10").Value = 1.73
$ws.Range("AL10").Value = 190

# Row 11
$ws.Range("Q11").Value = 1.64
$ws.Range("AO11").Value = 11.5

# Row 12
$ws.Range("I12").Value = 5.5
$ws.Range("T12").Value = 1.8
$ws.Range("AB12").Value = 10
$ws.Range("AD12").Value = 21
$ws.Range("AG12").Value = 11
$ws.Range("AH12").Value = 23
$ws.Range("AK12").Value = 26

# Row 13
$ws.Range("F13").Value = 3.5
$ws.Range("G13").Value = 3.9
$ws.Range("H13").Value = 2.14
$ws.Range("I13").Value = 2.26
$ws.Range("J13").Value = 3.5
$ws.Range("P13").Value = 1.97
$ws.Range("Q13").Value = 1.89

# Row 14
$ws.Range("P14").Value = 1.79
$ws.Range("T14").Value = 1.87
$ws.Range("AJ14").Value = 120

# Row 15
$ws.Range("F15").Value = 4.4
$ws.Range("G15").Value = 5.4
$ws.Range("H15").Value = 1.85
$ws.Range("I15").Value = 2.06
$ws.Range("J15").Value = 3.4
$ws.Range("P15").Value = 1.79
$ws.Range("Q15").Value = 2.02
$ws.Range("T15").Value = 1.86
$ws.Range("U15").Value = 1.92
$ws.Range("Y15").Value = 8.6
$ws.Range("AB15").Value = 17.5
$ws.Range("AM15").Value = 150
$ws.Range("AN15").Value = 100

# Row 16
$ws.Range("O16").Value = 1.26
$ws.Range("S16").Value = 3
$ws.Range("AH16").Value = 16.5

# Row 17
$ws.Range("F17").Value = 2.82

# Row 18
$ws.Range("G18").Value = 1.7
$ws.Range("J18").Value = 4.2

# Row 19
$ws.Range("F19").Value = 2.4
$ws.Range("G19").Value = 2.52
$ws.Range("I19").Value = 3.45
$ws.Range("M19").Value = 1.1
$ws.Range("T19").Value = 1.92
$ws.Range("U19").Value = 1.94
$ws.Range("X19").Value = 11
$ws.Range("Y19").Value = 11
$ws.Range("Z19").Value = 24
$ws.Range("AA19").Value = 80
$ws.Range("AB19").Value = 8.8
$ws.Range("AC19").Value = 7.8
$ws.Range("AD19").Value = 16
$ws.Range("AE19").Value = 55
$ws.Range("AF19").Value = 15.5
$ws.Range("AG19").Value = 12
$ws.Range("AH19").Value = 23
$ws.Range("AI19").Value = 70
$ws.Range("AJ19").Value = 42
$ws.Range("AK19").Value = 32
$ws.Range("AL19").Value = 60
$ws.Range("AM19").Value = 150
$ws.Range("AN19").Value = 29
$ws.Range("AO19").Value = 60

# Row 21
$ws.Range("N21").Value = 2.84
$ws.Range("Y21").Value = 10.5
$ws.Range("AF21").Value = 1000
